$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14809
$ws1.Range("F3").Value = 18310
$ws1.Range("F5").Value = 104
$ws1.Range("F9").Value = 53
$ws1.Range("F14").Value = 87
$ws1.Range("F17").Value = 1388
$ws1.Range("F19").Value = 79
$ws1.Range("F21").Value = 224
$ws1.Range("F22").Value = 7579
$ws1.Range("F24").Value = 16
$ws1.Range("F26").Value = 1204
$ws1.Range("F28").Value = 5923
$ws1.Range("F29").Value = 89
$ws1.Range("F30").Value = 57
$ws1.Range("F34").Value = 5252

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14809
$ws4.Range("F3").Value = 18310
$ws4.Range("F5").Value = 104
$ws4.Range("F9").Value = 53
$ws4.Range("F14").Value = 87
$ws4.Range("F17").Value = 1388
$ws4.Range("F19").Value = 79
$ws4.Range("F22").Value = 224
$ws4.Range("F23").Value = 7579
$ws4.Range("F25").Value = 16
$ws4.Range("F27").Value = 1204
$ws4.Range("F30").Value = 5923
$ws4.Range("F31").Value = 89
$ws4.Range("F32").Value = 57
$ws4.Range("F36").Value = 5252
